# Edit script: insert a "Status" column after "Mutual Fund" column, rename/shift the
# month columns (Jan_2026, Dec_2025, Oct_2025 replacing Nov_2025), recompute MoM/QoQ,
# and refresh the holdings rows (including several exits/new entries), per the
# "updated data from quant engine" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column D (Status) by shifting existing D:H one column to the right ---
$ws.Columns.Item(4).Insert()

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "ISIN"
$ws.Cells.Item(1, 2).Value = "Stock Name"
$ws.Cells.Item(1, 3).Value = "Mutual Fund"
$ws.Cells.Item(1, 4).Value = "Status"
$ws.Cells.Item(1, 5).Value = "Jan_2026"
$ws.Cells.Item(1, 6).Value = "Dec_2025"
$ws.Cells.Item(1, 7).Value = "Oct_2025"
$ws.Cells.Item(1, 8).Value = "MoM"
$ws.Cells.Item(1, 9).Value = "QoQ"

# Copy the header style (bold/centered/bordered) from an existing header cell (A1) onto
# the newly created "Status" header cell so it matches the rest of row 1.
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fund holdings data (rows 2-24) ---
# Columns: ISIN, Stock Name, MutualFund, Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ
$data = @(
    @("INE781S01027", "Ventive Hospitality Limited",                    "quant Consumption Fund", "Reducing Consistently", 9.802079000000001,  10.033282,           9.804278,   -0.2312029999999989,   -0.002198999999999174),
    @("INE180C01042", "Capri Global Capital Limited",                   "quant Consumption Fund", "Adding",                9.526871999999999,  9.512081999999999,  9.627039,    0.01478999999999964,  -0.1001670000000008),
    @("INE768C01028", "Zydus Wellness Ltd",                             "quant Consumption Fund", "Reducing Consistently", 8.403816000000001,  9.965963,            10.230047,  -1.562147,             -1.826231),
    @("INE016A01026", "Dabur India Limited",                            "quant Consumption Fund", "Adding Consistently",   7.315448,            3.949244,            0,           3.366204,              7.315448),
    @("INE917I01010", "Bajaj Auto Limited",                             "quant Consumption Fund", "Adding Consistently",   7.003441,            6.379967,            5.061371,    0.6234739999999999,    1.942069999999999),
    @("INE804L01022", "Medplus Health Services Limited",                "quant Consumption Fund", "Adding",                5.84704,             5.528791,            6.885661,    0.3182489999999998,   -1.038621),
    @("INE406A01037", "Aurobindo Pharma Limited",                       "quant Consumption Fund", "Adding Consistently",   5.578617,            3.282408,            2.880193,    2.296209,              2.698424000000001),
    @("INE04TZ01018", "ETHOS LIMITED",                                  "quant Consumption Fund", "Reducing",              5.232357,            5.71073,             0,          -0.4783729999999995,    5.232357),
    @("INE192A01025", "Tata Consumer Products Ltd",                     "quant Consumption Fund", "Adding Consistently",   3.671604,            0.493396,            0,           3.178208,              3.671604),
    @("INE179A01014", "Procter & Gamble Hygiene & Health Care Limited", "quant Consumption Fund", "Reducing",              3.426913,            3.535301,            3.316614,   -0.1083880000000002,    0.1102989999999999),
    @("INE01A001028", "Stanley Lifestyles Limited",                     "quant Consumption Fund", "Reducing Consistently", 1.90172,             3.168372,            5.623394,   -1.266652,             -3.721674),
    @("INE018E01016", "SBI Cards & Payment Services Ltd",               "quant Consumption Fund", "Reducing",              0.990131,            1.059535,            0.984699,   -0.06940399999999991,   0.005431999999999992),
    @("INE090A01021", "ICICI Bank Limited",                             "quant Consumption Fund", "Fresh Entry",           0.418628,            0,                   0,           0.418628,              0.418628),
    @("INE686F01025", "UNITED BREWERIES LIMITED",                       "quant Consumption Fund", "Complete Exit",         0,                   0,                   5.056928,    0,                    -5.056928),
    @("INE916U01025", "Sheela Foam Limited",                            "quant Consumption Fund", "Complete Exit",         0,                   0,                   2.955781,    0,                    -2.955781),
    @("INE854D01024", "United Spirits Limited",                        "quant Consumption Fund", "Complete Exit",         0,                   5.232451,            0,          -5.232451,              0),
    @("INE00H001014", "SWIGGY LIMITED",                                 "quant Consumption Fund", "Complete Exit",         0,                   0,                   2.641418,    0,                    -2.641418),
    @("INE669C01036", "Tech Mahindra Limited",                          "quant Consumption Fund", "Complete Exit",         0,                   0.272518,            0,          -0.272518,               0),
    @("INE484J01027", "Godrej Properties Limited",                      "quant Consumption Fund", "Complete Exit",         0,                   6.542531,            3.941633,   -6.542531,              -3.941633),
    @("INE364U01010", "Adani Green Energy Limited",                     "quant Consumption Fund", "Complete Exit",         0,                   0,                   3.123305,    0,                    -3.123305),
    @("INE196A01026", "Marico Limited",                                 "quant Consumption Fund", "Complete Exit",         0,                   6.171663,            0,          -6.171663,               0),
    @("INE14LE01019", "Aditya Birla Lifestyle Brands Limited",          "quant Consumption Fund", "Complete Exit",         0,                   0,                   5.125438,    0,                    -5.125438),
    @("INE202B01038", "Piramal Finance Ltd",                            "quant Consumption Fund", "Complete Exit",         0,                   0,                   2.132166,    0,                    -2.132166)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

$ws.Range("A1").Select() | Out-Null
